# Update countries & provincias Spain
# Applies the 27-May-2020 17:35 data refresh to the Pais sheet:
#  - swap the ranking of Chile and Arabia Saudita (Chile moved above
#    Arabia Saudita in row order, taking row 18, with Arabia Saudita
#    dropping to row 19 keeping its own stats)
#  - refresh the numeric stats (Casos totales, Nuevos casos, Casos
#    activos, Recuperados, Casos criticos, Muertes hoy, Muertes) for the
#    affected countries
#  - bump the "Datos actualizados" timestamp footer from 17:05 to 17:35

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Footer timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 27 de Mayo de 2020 a las 17:35"

# --- Row 4 : Estados Unidos -------------------------------------------
$ws.Range("B4").Value = 1730685
$ws.Range("C4").Value = 5410
$ws.Range("E4").Value = 1149539
$ws.Range("G4").Value = 253
$ws.Range("H4").Value = 100825

# --- Row 5 : Brasil -----------------------------------------------------
$ws.Range("E5").Value = 211314
$ws.Range("G5").Value = 51
$ws.Range("H5").Value = 24600

# --- Row 8 : Reino Unido -------------------------------------------------
$ws.Range("B8").Value = 267240
$ws.Range("C8").Value = 2013
$ws.Range("G8").Value = 412
$ws.Range("H8").Value = 37460

# --- Row 13 : Turquia -----------------------------------------------------
$ws.Range("B13").Value = 154820
$ws.Range("C13").Value = 4027
$ws.Range("D13").Value = 65944
$ws.Range("E13").Value = 84470
$ws.Range("G13").Value = 62
$ws.Range("H13").Value = 4406

# --- Row 16 : Canada -----------------------------------------------------
$ws.Range("B16").Value = 86939
$ws.Range("C16").Value = 292
$ws.Range("D16").Value = 45753
$ws.Range("E16").Value = 34515
$ws.Range("G16").Value = 32
$ws.Range("H16").Value = 6671

# --- Row 18/19 : Chile overtakes Arabia Saudita ---------------------------
# Row 18 now holds Chile's refreshed stats ...
$ws.Range("A18").Value = "Chile"
$ws.Range("B18").Value = 82289
$ws.Range("C18").Value = 4328
$ws.Range("D18").Value = 33540
$ws.Range("E18").Value = 47908
$ws.Range("G18").Value = 35
$ws.Range("H18").Value = 841

# ... and row 19 keeps Arabia Saudita's (previous row-18) stats
$ws.Range("A19").Value = "Arabia Saudita"
$ws.Range("B19").Value = 78541
$ws.Range("C19").Value = 1815
$ws.Range("D19").Value = 51022
$ws.Range("E19").Value = 27094
$ws.Range("G19").Value = 14
$ws.Range("H19").Value = 425

# --- Row 29 : Singapur -----------------------------------------------------
$ws.Range("D29").Value = 17276
$ws.Range("E29").Value = 15577

# --- Row 45 : Republica Dominicana ------------------------------------------
$ws.Range("B45").Value = 15723
$ws.Range("C45").Value = 459
$ws.Range("D45").Value = 8790
$ws.Range("E45").Value = 6459
$ws.Range("G45").Value = 6
$ws.Range("H45").Value = 474

# --- Row 81 ------------------------------------------------------------------
$ws.Range("B81").Value = 2903
$ws.Range("C81").Value = 11
$ws.Range("E81").Value = 1356

# --- Row 91 ------------------------------------------------------------------
$ws.Range("B91").Value = 1974
$ws.Range("C91").Value = 11
$ws.Range("D91").Value = 1724
$ws.Range("E91").Value = 168

# --- Row 103 -----------------------------------------------------------------
$ws.Range("B103").Value = 1453
$ws.Range("C103").Value = 134
$ws.Range("E103").Value = 711
